# #CRM-1168 Add brand column in Partner panel - Pending Spares - Download file
#
# Adds a new "Brand" column (column O) to the Spare Requested Parts export
# template:
#   O1 -> header text "Brand"            (same look as the other header cells)
#   O2 -> merge placeholder "{spare:brands}" (same look as the other placeholders)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing header style (bold font + grey fill + centered, the same
# formatting already used by A1:N1) for the new header cell instead of minting
# a brand new cell style.
$ws.Range("A1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Brand"

# Reuse the existing placeholder-row style (the same formatting already used
# by N2) for the new placeholder cell.
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("O2").Value = "{spare:brands}"

# Match the saved selection state of the edited workbook.
$ws.Range("N10").Select()
